$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.118.64"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.881.83"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.17%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.51"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("E6").Value = "  +0.12%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5038"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.29%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3837"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.55%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.08561"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -7.23%  "
$ws.Range("E10").Value = "  -1.15%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "41.87"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "6.269"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.92%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.878.81"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.81%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "20.58"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.17%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.218"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("E17").Value = "  -0.96%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "91.26"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.30%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06662"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.18%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "18.09"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  +0.18%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.101"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.75%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.157.93"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.20"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.19%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.265"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.41%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.597"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.73%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.095.69"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "20.74"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.57%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "156.38"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.28%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "126.39"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.45%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.1052"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.60%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.058"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  +0.66%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.610"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "9.727"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.79%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02463"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.59%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.06570"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2178"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("E39").Value = "  +0.32%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.6526"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.82%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.243"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -8.13%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "11.39"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "4.920"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.01%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.6202"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.59%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.18"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("E46").Value = "  -0.24%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.686"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.19%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.029"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("E49").Value = "  +1.46%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "121.07"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.96%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "80.88"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.26%  "
